$d = $word.ActiveDocument

# --- Text corrections (typo fixes), applied in left-to-right order ---
$fixes = @(
    @("Kilonrinnella,", "Kilonrinteessä,"),
    @("K-ruoka", "K-kauppa"),
    @("rautatieasema.", "rautatieasemaa."),
    @("olohuoneetta, muutta", "olohuonetta, mutta"),
    @("makuuhuoneetta.", "makuuhuonetta."),
    @("nojapuolia.", "tuolia."),
    @("työpuoli", "työtuoli"),
    @("ikkuna.", "ikkunaa."),
    @("Minullä", "Minulla"),
    @("naapuria", "naapuri"),
    @("ystävällinen.", "ystävälliset.")
)

foreach ($fix in $fixes) {
    $old = $fix[0]
    $new = $fix[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# --- Paragraph-level formatting ---
$para = $d.Paragraphs(1)
$para.Style = "Normal"
$para.Format.SpaceBefore = 0
$para.Format.SpaceAfter = 8

# --- Run-level formatting: bump font size to 14pt (sz/szCs = 28) ---
# Apply across the whole paragraph (including its end-of-paragraph mark)
# so both the runs and the paragraph mark pick up the new size.
$d.Paragraphs(1).Range.Font.Size = 14
$d.Paragraphs(1).Range.Font.SizeBi = 14

# --- Section/page setup updates ---
$ps = $d.Sections(1).PageSetup
$ps.HeaderDistance = 0
$ps.FooterDistance = 0
